$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 0.009
$ws.Range("E3").Value = 0.225
$ws.Range("F3").Value = 0.002
$ws.Range("G3").Value = 0.006
$ws.Range("I3").Value = 0.004

$ws.Range("D4").Value = 0.182
$ws.Range("E4").Value = 0.395
$ws.Range("F4").Value = 0.154
$ws.Range("G4").Value = 0.202
$ws.Range("H4").Value = 0.111
$ws.Range("I4").Value = 0.192

$ws.Range("D5").Value = 0.81
$ws.Range("E5").Value = 0.591
$ws.Range("F5").Value = 0.842
$ws.Range("G5").Value = 0.814
$ws.Range("H5").Value = 0.87
$ws.Range("I5").Value = 0.812

$ws.Range("D6").Value = 0.998
$ws.Range("E6").Value = 0.777
$ws.Range("F6").Value = 0.999
$ws.Range("G6").Value = 0.996
$ws.Range("H6").Value = 0.999
$ws.Range("I6").Value = 0.997

$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 0.908
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 1
$ws.Range("I7").Value = 1

$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 0.972
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1
$ws.Range("H8").Value = 1
$ws.Range("I8").Value = 1

$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 0.995
$ws.Range("G9").Value = 1
$ws.Range("I9").Value = 1

$ws.Range("E10").Value = 0.999

$ws.Range("E11").Value = 1

$ws.Range("D15").Value = 0.978
$ws.Range("E15").Value = 0.809
$ws.Range("F15").Value = 0.955
$ws.Range("G15").Value = 0.97
$ws.Range("H15").Value = 0.92
$ws.Range("I15").Value = 0.977

$ws.Range("D16").Value = 0.044
$ws.Range("E16").Value = 0.188
$ws.Range("F16").Value = 0.079
$ws.Range("G16").Value = 0.059
$ws.Range("H16").Value = 0.121
$ws.Range("I16").Value = 0.041

$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 0.014
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0

$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 0.001
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0

$ws.Range("D19").Value = 0
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 0
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0

$ws.Range("E20").Value = 0

$ws.Range("D25").Value = 0.001
$ws.Range("E25").Value = 0.028
$ws.Range("F25").Value = 0.056
$ws.Range("G25").Value = 0.003
$ws.Range("H25").Value = 0.365

$ws.Range("D26").Value = 0.167
$ws.Range("E26").Value = 0.291
$ws.Range("F26").Value = 0.296
$ws.Range("G26").Value = 0.177
$ws.Range("H26").Value = 0.44
$ws.Range("I26").Value = 0.169

$ws.Range("D27").Value = 0.81
$ws.Range("E27").Value = 0.716
$ws.Range("F27").Value = 0.694
$ws.Range("G27").Value = 0.817
$ws.Range("H27").Value = 0.532
$ws.Range("I27").Value = 0.818

$ws.Range("D28").Value = 0.994
$ws.Range("E28").Value = 0.96
$ws.Range("F28").Value = 0.953
$ws.Range("G28").Value = 0.997
$ws.Range("H28").Value = 0.603
$ws.Range("I28").Value = 0.999

$ws.Range("D29").Value = 1
$ws.Range("E29").Value = 0.999
$ws.Range("F29").Value = 0.996
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 0.689
$ws.Range("I29").Value = 1

$ws.Range("D30").Value = 1
$ws.Range("E30").Value = 1
$ws.Range("F30").Value = 1
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 0.781
$ws.Range("I30").Value = 1

$ws.Range("D31").Value = 1
$ws.Range("E31").Value = 1
$ws.Range("F31").Value = 1
$ws.Range("G31").Value = 1
$ws.Range("H31").Value = 0.838
$ws.Range("I31").Value = 1

$ws.Range("H32").Value = 0.88

$ws.Range("H33").Value = 0.919
